# Add demos corresponding to the last 10 CO-published papers to the
# "v2.5" sheet of paper_demo.xlsx (rows 3-11, columns A (paper), and the
# Figure 1..10 demo columns B..K as applicable).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("v2.5")

# --- Pass 1: paper/article names (column A), top to bottom ---
$ws.Range("A3").Value  = "Li et al. (2012)"
$ws.Range("A4").Value  = "Buck et al. (2012)"
$ws.Range("A5").Value  = "Murphy (2012)"
$ws.Range("A6").Value  = "Peng and Murphy (2011)"
$ws.Range("A7").Value  = "Shariff et al (2011)"
$ws.Range("A8").Value  = "Murphy (2010)"
$ws.Range("A9").Value  = "Shariff et al. (2010)"
$ws.Range("A10").Value = "Peng et al. (2009)"
$ws.Range("A11").Value = "Rohde et al. (2008)"

# --- Pass 2: demo references for each paper's figures ---
$ws.Range("E3").Value = "demo3D04"
$ws.Range("C3").Value = "demo3D00, demo3D02"
$ws.Range("F3").Value = "demo3Dimg2microtubule_model"

$ws.Range("B4").Value = "demo2D01"
$ws.Range("C4").Value = "demo3D11"
$ws.Range("E4").Value = "demo2D02, demo3D09"
$ws.Range("G4").Value = "demo3D04"

$ws.Range("B5").Value = "demo2D01"
$ws.Range("C5").Value = "demo3D11"

$ws.Range("B6").Value = "demo3D11"
$ws.Range("I6").Value = "demo3D02"

$ws.Range("B7").Value = "demo3D04"
$ws.Range("C7").Value = "demo3D04"

$ws.Range("E8").Value = "demo3D01"
$ws.Range("G8").Value = "demo2D00"

$ws.Range("B9").Value = "demo3D01, demo3D14"
$ws.Range("E9").Value = "demo3D07, demo3D06"
$ws.Range("G9").Value = "demo3Dimg2microtubule_model"
$ws.Range("H9").Value = "demo3Dimg2microtubule_model"
$ws.Range("I9").Value = "demo3Dimg2microtubule_model"

$ws.Range("C10").Value = "demo2D04, demo3D15, demo3D20"
$ws.Range("F10").Value = "demo3D15"

# Page setup: force portrait orientation (as saved by the newer Excel build).
$ws.PageSetup.Orientation = 1

# Final cursor/selection position left by the editing session.
$ws.Range("E13").Select() | Out-Null
